$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2: Mahindra Tractors | 275 DI TU XP Plus | [images...]
$ws.Range("B2").Value2 = "275 DI TU XP Plus"
$ws.Range("C2").Value2 = "['275 DI TU XP Plusimg0-275-di-tu-xp-plus-1632304804.png', '275 DI TU XP Plusimg1-275-di-tu-xp-plus-1632304804.png', '275 DI TU XP Plusimg2-mqdefault.png']"

# Update row 3: Mahindra Tractors | Arjun 555 DI | [images...]
$ws.Range("A3").Value2 = "Mahindra Tractors"
$ws.Range("B3").Value2 = "Arjun 555 DI"
$ws.Range("C3").Value2 = "['Arjun 555 DIimg0-arjun-555-di-1632207634.png', 'Arjun 555 DIimg1-mqdefault.png', 'Arjun 555 DIimg2-mqdefault.png', 'Arjun 555 DIimg3-arjun-555-di-1632207634.png']"

# Remove old rows 4-6 (John Deere / Swaraj entries no longer present)
$ws.Range("A4:C6").Delete() | Out-Null
